$d = $word.ActiveDocument

# 1. Update the main title
$d.Content.Find.Execute(
    "Play Golden Ark Free: Review of Egyptian-Themed Slot Game", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Golden Ark for Free - Exciting Egyptian-Themed Slot Game", 2) | Out-Null

# 2. Turn the "Meta description" paragraph into the new bold heading that will be
#    relocated near the end of the document (right after the "What we don't like" list).
#    First strip the trailing descriptive sentence, leaving just the bold run.
$d.Content.Find.Execute(
    ": Discover the exciting features of Golden Ark, an Egyptian-themed slot game with free spins and gamble feature. Play for free and read our unbiased review.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Rename the remaining bold "Meta description" run to the new heading text.
$d.Content.Find.Execute(
    "Meta description", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Golden Ark for Free - Exciting Egyptian-Themed Slot Game", 2) | Out-Null

# Cut the (now repurposed) paragraph out of its original position (right under the title).
$metaParaRange = $d.Paragraphs(2).Range
$metaParaRange.Cut() | Out-Null

# 3. Update the "What we like" bullet points.
$d.Content.Find.Execute(
    "Exciting Gamble feature", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting Gamble feature for more chances to win", 2) | Out-Null

$d.Content.Find.Execute(
    "Free spins triggered by Book of Ra symbol", $true, $false, $false, $false, $false,
    $true, 1, $false, "Free spins feature can be triggered for up to 10 free spins", 2) | Out-Null

$d.Content.Find.Execute(
    "Autoplay function available", $true, $false, $false, $false, $false,
    $true, 1, $false, "Autoplay function for convenient gameplay", 2) | Out-Null

$d.Content.Find.Execute(
    "Visually appealing Egyptian theme", $true, $false, $false, $false, $false,
    $true, 1, $false, "Similar games with ancient civilization theme available", 2) | Out-Null

# 4. Update the "What we don't like" bullet points.
$d.Content.Find.Execute(
    "Limited paylines", $true, $false, $false, $false, $false,
    $true, 1, $false, "Risk of losing winnings with the Gamble feature", 2) | Out-Null

$d.Content.Find.Execute(
    "No progressive jackpot", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited number of paylines compared to some other slot games", 2) | Out-Null

# 5. Paste the relocated heading paragraph right after the last "don't like" bullet.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Limited number of paylines compared to some other slot games*") {
        $targetIdx = $i
    }
}
$targetPara = $d.Paragraphs($targetIdx)
$insertPoint = $d.Range($targetPara.Range.End, $targetPara.Range.End)
$insertPoint.Paste()

# 6. Replace the closing image-prompt paragraph with the new promotional sentence
#    (keeping the existing italic formatting of that run).
$d.Content.Find.Execute(
    "Prompt: Create a feature image for Golden Ark that captures the thrilling adventure of exploring an ancient pyramid in search of treasure. The image should be in cartoon style and prominently feature a happy Maya warrior with glasses. The warrior should be depicted holding a book or treasure chest, with hieroglyphics and ancient artifacts surrounding him. The colors should be vibrant and eye-catching, with a sense of mystery and excitement. Use your creativity to bring the world of Golden Ark to life and entice players to embark on this journey of discovery.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Golden Ark for free and experience the thrill of ancient Egypt in this exciting slot game.", 2) | Out-Null
